# Updates cryptocurrency price/volume/hour data per the commit
# "Updated symbol list on Tue Feb 14 12:22:41 UTC 2023 with GitHub Actions".
# Each target cell is written as literal text (not auto-converted to a
# number/date by Excel) by building it via a `="..."` text formula and then
# collapsing that formula down to its static value with Copy + PasteSpecial
# (paste values only), which avoids introducing any new/changed cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $cell = $ws.Range($range)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue 'D2' '293.94'
Set-TextValue 'E2' '-0.31%'
Set-TextValue 'G2' '12'
Set-TextValue 'D3' '40.54'
Set-TextValue 'E3' '0.59%'
Set-TextValue 'G3' '12'
Set-TextValue 'D4' '5.014'
Set-TextValue 'E4' '0.10%'
Set-TextValue 'G4' '12'
Set-TextValue 'D5' '0.07410'
Set-TextValue 'E5' '0.21%'
Set-TextValue 'G5' '12'
Set-TextValue 'D6' '1.575'
Set-TextValue 'E6' '-0.53%'
Set-TextValue 'G6' '12'
Set-TextValue 'D7' '0.9254'
Set-TextValue 'E7' '0.50%'
Set-TextValue 'G7' '12'
Set-TextValue 'D8' '2.350'
Set-TextValue 'E8' '-2.04%'
Set-TextValue 'G8' '12'
Set-TextValue 'D9' '0.1206'
Set-TextValue 'E9' '1.34%'
Set-TextValue 'G9' '12'
Set-TextValue 'E10' '3.42%'
Set-TextValue 'G10' '12'
Set-TextValue 'D11' '0.04383'
Set-TextValue 'E11' '5.33%'
Set-TextValue 'G11' '12'
Set-TextValue 'D12' '0.08768'
Set-TextValue 'E12' '1.12%'
Set-TextValue 'G12' '12'
Set-TextValue 'D13' '0.1055'
Set-TextValue 'E13' '0.14%'
Set-TextValue 'G13' '12'
Set-TextValue 'D14' '0.001278'
Set-TextValue 'E14' '0.38%'
Set-TextValue 'G14' '12'
Set-TextValue 'D15' '0.006028'
Set-TextValue 'E15' '3.88%'
Set-TextValue 'G15' '12'
Set-TextValue 'E16' '-2.11%'
Set-TextValue 'G16' '12'
Set-TextValue 'D17' '4.295'
Set-TextValue 'E17' '-0.08%'
Set-TextValue 'G17' '12'
Set-TextValue 'E18' '0.66%'
Set-TextValue 'G18' '12'
Set-TextValue 'D19' '7.852'
Set-TextValue 'E19' '3.84%'
Set-TextValue 'G19' '12'
Set-TextValue 'D20' '0.1391'
Set-TextValue 'E20' '3.64%'
Set-TextValue 'G20' '12'
Set-TextValue 'D21' '0.2883'
Set-TextValue 'E21' '2.84%'
Set-TextValue 'G21' '12'
Set-TextValue 'D22' '0.03933'
Set-TextValue 'E22' '2.43%'
Set-TextValue 'G22' '12'
Set-TextValue 'D23' '0.001260'
Set-TextValue 'E23' '-1.61%'
Set-TextValue 'G23' '12'
Set-TextValue 'D24' '0.003784'
Set-TextValue 'E24' '-3.10%'
Set-TextValue 'G24' '12'
Set-TextValue 'D25' '0.0001230'
Set-TextValue 'E25' '-4.77%'
Set-TextValue 'G25' '12'
Set-TextValue 'D26' '0.0003726'
Set-TextValue 'E26' '-0.03%'
Set-TextValue 'G26' '12'
Set-TextValue 'G27' '12'
Set-TextValue 'G28' '12'
Set-TextValue 'G29' '12'
Set-TextValue 'G30' '12'
Set-TextValue 'G31' '12'
Set-TextValue 'G32' '12'
Set-TextValue 'G33' '12'
Set-TextValue 'G34' '12'
Set-TextValue 'G35' '12'
Set-TextValue 'G36' '12'
Set-TextValue 'G37' '12'
Set-TextValue 'D38' '0.02327'
Set-TextValue 'E38' '-0.03%'
Set-TextValue 'G38' '12'
Set-TextValue 'D39' '0.05089'
Set-TextValue 'E39' '1.05%'
Set-TextValue 'G39' '12'
Set-TextValue 'D40' '0.006318'
Set-TextValue 'E40' '40.23%'
Set-TextValue 'G40' '12'
Set-TextValue 'D41' '0.007806'
Set-TextValue 'E41' '1.33%'
Set-TextValue 'G41' '12'
Set-TextValue 'D42' '0.1292'
Set-TextValue 'E42' '1.20%'
Set-TextValue 'G42' '12'
Set-TextValue 'D43' '0.007391'
Set-TextValue 'E43' '-0.16%'
Set-TextValue 'G43' '12'
Set-TextValue 'D44' '0.007286'
Set-TextValue 'E44' '4.84%'
Set-TextValue 'G44' '12'
Set-TextValue 'D45' '0.2919'
Set-TextValue 'E45' '-9.12%'
Set-TextValue 'G45' '12'
Set-TextValue 'D46' '0.00006213'
Set-TextValue 'E46' '-3.90%'
Set-TextValue 'G46' '12'
Set-TextValue 'E47' '0.00%'
Set-TextValue 'G47' '12'
Set-TextValue 'D48' '0.04682'
Set-TextValue 'E48' '-81.41%'
Set-TextValue 'G48' '12'
Set-TextValue 'D49' '0.004204'
Set-TextValue 'E49' '-0.02%'
Set-TextValue 'G49' '12'
Set-TextValue 'E50' '0.00%'
Set-TextValue 'G50' '12'
Set-TextValue 'E51' '0.00%'
Set-TextValue 'G51' '12'

$excel.CutCopyMode = $false
